$d = $word.ActiveDocument
$ir = $d.Content
$ir.Collapse(0)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>Awards/Achievements</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
    <w:t>Eleven Fifty Core Value Award: Teamwork &#8211; December 2021</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NoSpacing"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:b/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cstheme="minorHAnsi"/>
      <w:bCs/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
    <w:t>Voted by my classmates as the individual who most exemplified a commitment to working together as one in utilizing the team&#8217;s combined strengths to achieve a shared vision.</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$ir.InsertXML($xml)
Write-Output "Done. Para count: $($d.Paragraphs.Count)"
